$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.072208172810626
$ws.Cells.Item(2, 4).Value = 1.072030756521074
$ws.Cells.Item(2, 5).Value = 1.075943553138918
$ws.Cells.Item(2, 6).Value = 1.085447045681993
$ws.Cells.Item(2, 9).Value = 1.045430307543311
$ws.Cells.Item(2, 10).Value = 1.077128879515603
$ws.Cells.Item(2, 11).Value = 1.074726353444159
$ws.Cells.Item(2, 12).Value = 1.078628789076833
$ws.Cells.Item(2, 13).Value = 1.088107452401265
$ws.Cells.Item(2, 14).Value = 1.078658526072924
$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.073897414423404
$ws.Cells.Item(3, 4).Value = 1.073362654876628
$ws.Cells.Item(3, 5).Value = 1.077433414607287
$ws.Cells.Item(3, 6).Value = 1.086982225004268
$ws.Cells.Item(3, 9).Value = 1.045811968428813
$ws.Cells.Item(3, 10).Value = 1.078472878139405
$ws.Cells.Item(3, 11).Value = 1.075873647800539
$ws.Cells.Item(3, 12).Value = 1.079934415942693
$ws.Cells.Item(3, 13).Value = 1.089460104564123
$ws.Cells.Item(3, 14).Value = 1.080004433328931
$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.07498858742434
$ws.Cells.Item(4, 4).Value = 1.07422256934854
$ws.Cells.Item(4, 5).Value = 1.078395897602642
$ws.Cells.Item(4, 6).Value = 1.087974049717594
$ws.Cells.Item(4, 9).Value = 1.046056705516309
$ws.Cells.Item(4, 10).Value = 1.079340287944792
$ws.Cells.Item(4, 11).Value = 1.076613566823521
$ws.Cells.Item(4, 12).Value = 1.080777172615589
$ws.Cells.Item(4, 13).Value = 1.09033331298549
$ws.Cells.Item(4, 14).Value = 1.080873074955736
$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.075446877589647
$ws.Cells.Item(5, 4).Value = 1.074583627166517
$ws.Cells.Item(5, 5).Value = 1.078800161393003
$ws.Cells.Item(5, 6).Value = 1.088390653222823
$ws.Cells.Item(5, 9).Value = 1.04615906358958
$ws.Cells.Item(5, 10).Value = 1.079704417546836
$ws.Cells.Item(5, 11).Value = 1.07692404775108
$ws.Cells.Item(5, 12).Value = 1.081130979719386
$ws.Cells.Item(5, 13).Value = 1.090699927685366
$ws.Cells.Item(5, 14).Value = 1.081237721663581
$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.075523801100411
$ws.Cells.Item(6, 4).Value = 1.074644224184347
$ws.Cells.Item(6, 5).Value = 1.078868017964874
$ws.Cells.Item(6, 6).Value = 1.088460581952467
$ws.Cells.Item(6, 9).Value = 1.046176218997199
$ws.Cells.Item(6, 10).Value = 1.07976552570611
$ws.Cells.Item(6, 11).Value = 1.076976145001797
$ws.Cells.Item(6, 12).Value = 1.081190357116916
$ws.Cells.Item(6, 13).Value = 1.090761455871181
$ws.Cells.Item(6, 14).Value = 1.08129891660345
$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.074994712829632
$ws.Cells.Item(7, 4).Value = 1.07422739557909
$ws.Cells.Item(7, 5).Value = 1.078401300813604
$ws.Cells.Item(7, 6).Value = 1.087979617794718
$ws.Cells.Item(7, 9).Value = 1.046058075306129
$ws.Cells.Item(7, 10).Value = 1.079345155531027
$ws.Cells.Item(7, 11).Value = 1.076617717757854
$ws.Cells.Item(7, 12).Value = 1.08078190210894
$ws.Cells.Item(7, 13).Value = 1.090338213592071
$ws.Cells.Item(7, 14).Value = 1.080877949454501
$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.072779455644645
$ws.Cells.Item(8, 4).Value = 1.072481278525498
$ws.Cells.Item(8, 5).Value = 1.076447386431177
$ws.Cells.Item(8, 6).Value = 1.085966190704515
$ws.Cells.Item(8, 9).Value = 1.04555975361005
$ws.Cells.Item(8, 10).Value = 1.077583560567763
$ws.Cells.Item(8, 11).Value = 1.075114600069533
$ws.Cells.Item(8, 12).Value = 1.079070465510181
$ws.Cells.Item(8, 13).Value = 1.088565016514029
$ws.Cells.Item(8, 14).Value = 1.079113852824332
$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.068861004459297
$ws.Cells.Item(9, 4).Value = 1.069389400120192
$ws.Cells.Item(9, 5).Value = 1.072992023924897
$ws.Cells.Item(9, 6).Value = 1.082406091072583
$ws.Cells.Item(9, 9).Value = 1.044664495167863
$ws.Cells.Item(9, 10).Value = 1.074461812029879
$ws.Cells.Item(9, 11).Value = 1.072446770963086
$ws.Cells.Item(9, 12).Value = 1.076038470150415
$ws.Cells.Item(9, 13).Value = 1.085424363913544
$ws.Cells.Item(9, 14).Value = 1.075987671045454
$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.066238005587018
$ws.Cells.Item(10, 4).Value = 1.06731757653338
$ws.Cells.Item(10, 5).Value = 1.07067961924238
$ws.Cells.Item(10, 6).Value = 1.080023941424892
$ws.Cells.Item(10, 9).Value = 1.04405594622136
$ws.Cells.Item(10, 10).Value = 1.072368280377774
$ws.Cells.Item(10, 11).Value = 1.0706548901889
$ws.Cells.Item(10, 12).Value = 1.074005721658344
$ws.Cells.Item(10, 13).Value = 1.083319285692774
$ws.Cells.Item(10, 14).Value = 1.073891166338272
$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.065099526407485
$ws.Cells.Item(11, 4).Value = 1.066417837230679
$ws.Cells.Item(11, 5).Value = 1.069676105221957
$ws.Cells.Item(11, 6).Value = 1.078990243979441
$ws.Cells.Item(11, 9).Value = 1.043789620291422
$ws.Cells.Item(11, 10).Value = 1.071458704713226
$ws.Cells.Item(11, 11).Value = 1.069875722366402
$ws.Cells.Item(11, 12).Value = 1.073122696825941
$ws.Cells.Item(11, 13).Value = 1.08240496526894
$ws.Cells.Item(11, 14).Value = 1.072980298971951
$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.064676225033307
$ws.Cells.Item(12, 4).Value = 1.06608322994097
$ws.Cells.Item(12, 5).Value = 1.069303010092781
$ws.Cells.Item(12, 6).Value = 1.078605939700836
$ws.Cells.Item(12, 9).Value = 1.043690267758681
$ws.Cells.Item(12, 10).Value = 1.071120377231516
$ws.Cells.Item(12, 11).Value = 1.069585804301572
$ws.Cells.Item(12, 12).Value = 1.072794266801465
$ws.Cells.Item(12, 13).Value = 1.082064914069168
$ws.Cells.Item(12, 14).Value = 1.072641491026409
$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.064767043814888
$ws.Cells.Item(13, 4).Value = 1.066155022778515
$ws.Cells.Item(13, 5).Value = 1.06938305606439
$ws.Cells.Item(13, 6).Value = 1.078688389970978
$ws.Cells.Item(13, 9).Value = 1.043711598603864
$ws.Cells.Item(13, 10).Value = 1.071192971088475
$ws.Cells.Item(13, 11).Value = 1.069648015500832
$ws.Cells.Item(13, 12).Value = 1.072864736014272
$ws.Cells.Item(13, 13).Value = 1.082137875912119
$ws.Cells.Item(13, 14).Value = 1.072714187974972
$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.065064544819603
$ws.Cells.Item(14, 4).Value = 1.066390186785672
$ws.Cells.Item(14, 5).Value = 1.069645272169161
$ws.Cells.Item(14, 6).Value = 1.078958484328586
$ws.Cells.Item(14, 9).Value = 1.043781416515792
$ws.Cells.Item(14, 10).Value = 1.071430748118574
$ws.Cells.Item(14, 11).Value = 1.069851767909506
$ws.Cells.Item(14, 12).Value = 1.073095557635009
$ws.Cells.Item(14, 13).Value = 1.082376865399553
$ws.Cells.Item(14, 14).Value = 1.072952302675727
$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.065247789015466
$ws.Cells.Item(15, 4).Value = 1.06653502526466
$ws.Cells.Item(15, 5).Value = 1.069806786083399
$ws.Cells.Item(15, 6).Value = 1.079124852596335
$ws.Cells.Item(15, 9).Value = 1.043824376918416
$ws.Cells.Item(15, 10).Value = 1.071577187723506
$ws.Cells.Item(15, 11).Value = 1.069977239911331
$ws.Cells.Item(15, 12).Value = 1.073237716476047
$ws.Cells.Item(15, 13).Value = 1.082524057178846
$ws.Cells.Item(15, 14).Value = 1.073098950241694
$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.066313504605817
$ws.Cells.Item(16, 4).Value = 1.0673772331298
$ws.Cells.Item(16, 5).Value = 1.070746171322119
$ws.Cells.Item(16, 6).Value = 1.080092497014
$ws.Cells.Item(16, 9).Value = 1.044073561677926
$ws.Cells.Item(16, 10).Value = 1.072428580502587
$ws.Cells.Item(16, 11).Value = 1.070706531302491
$ws.Cells.Item(16, 12).Value = 1.074064264611757
$ws.Cells.Item(16, 13).Value = 1.08337990611691
$ws.Cells.Item(16, 14).Value = 1.073951552096178
$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.066981266674723
$ws.Cells.Item(17, 4).Value = 1.067904818119008
$ws.Cells.Item(17, 5).Value = 1.071334818982005
$ws.Cells.Item(17, 6).Value = 1.08069887523212
$ws.Cells.Item(17, 9).Value = 1.044229111168384
$ws.Cells.Item(17, 10).Value = 1.072961809393014
$ws.Cells.Item(17, 11).Value = 1.071163114348186
$ws.Cells.Item(17, 12).Value = 1.074581971426811
$ws.Cells.Item(17, 13).Value = 1.08391599877735
$ws.Cells.Item(17, 14).Value = 1.074485538232792
$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.067370501066232
$ws.Cells.Item(18, 4).Value = 1.068212296792879
$ws.Cells.Item(18, 5).Value = 1.071677953005001
$ws.Cells.Item(18, 6).Value = 1.081052352803943
$ws.Cells.Item(18, 9).Value = 1.044319568683802
$ws.Cells.Item(18, 10).Value = 1.073272537993952
$ws.Cells.Item(18, 11).Value = 1.071429116456993
$ws.Cells.Item(18, 12).Value = 1.074883668546218
$ws.Cells.Item(18, 13).Value = 1.084228422342516
$ws.Cells.Item(18, 14).Value = 1.074796708103984
$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.067503176249566
$ws.Cells.Item(19, 4).Value = 1.068317096494784
$ws.Cells.Item(19, 5).Value = 1.071794916829174
$ws.Cells.Item(19, 6).Value = 1.081172843790049
$ws.Cells.Item(19, 9).Value = 1.044350366350737
$ws.Cells.Item(19, 10).Value = 1.073378438730156
$ws.Cells.Item(19, 11).Value = 1.071519763131218
$ws.Cells.Item(19, 12).Value = 1.074986493543762
$ws.Cells.Item(19, 13).Value = 1.084334905174438
$ws.Cells.Item(19, 14).Value = 1.07490275923138
$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.066909649106367
$ws.Cells.Item(20, 4).Value = 1.067848239452039
$ws.Cells.Item(20, 5).Value = 1.071271684891757
$ws.Cells.Item(20, 6).Value = 1.080633838655089
$ws.Cells.Item(20, 9).Value = 1.044212450326483
$ws.Cells.Item(20, 10).Value = 1.07290462954656
$ws.Cells.Item(20, 11).Value = 1.071114159960786
$ws.Cells.Item(20, 12).Value = 1.074526454643452
$ws.Cells.Item(20, 13).Value = 1.0838585091307
$ws.Cells.Item(20, 14).Value = 1.074428277184398
$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.064976949877381
$ws.Cells.Item(21, 4).Value = 1.066320948077607
$ws.Cells.Item(21, 5).Value = 1.069568065681215
$ws.Cells.Item(21, 6).Value = 1.078878957828961
$ws.Cells.Item(21, 9).Value = 1.043760868702801
$ws.Cells.Item(21, 10).Value = 1.071360741784421
$ws.Cells.Item(21, 11).Value = 1.069791781790275
$ws.Cells.Item(21, 12).Value = 1.073027598516831
$ws.Cells.Item(21, 13).Value = 1.082306500955577
$ws.Cells.Item(21, 14).Value = 1.072882196924549
$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.063759351263831
$ws.Cells.Item(22, 4).Value = 1.065358336870543
$ws.Cells.Item(22, 5).Value = 1.068494929271717
$ws.Cells.Item(22, 6).Value = 1.077773604532852
$ws.Cells.Item(22, 9).Value = 1.043474467857289
$ws.Cells.Item(22, 10).Value = 1.07038731038635
$ws.Cells.Item(22, 11).Value = 1.068957450260523
$ws.Cells.Item(22, 12).Value = 1.072082684779306
$ws.Cells.Item(22, 13).Value = 1.081328188013989
$ws.Cells.Item(22, 14).Value = 1.07190738314224
$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.064405058525692
$ws.Cells.Item(23, 4).Value = 1.06586886052898
$ws.Cells.Item(23, 5).Value = 1.069064012675969
$ws.Cells.Item(23, 6).Value = 1.078359765514451
$ws.Cells.Item(23, 9).Value = 1.043626530012234
$ws.Cells.Item(23, 10).Value = 1.07090360687511
$ws.Cells.Item(23, 11).Value = 1.06940002306568
$ws.Cells.Item(23, 12).Value = 1.072583843951574
$ws.Cells.Item(23, 13).Value = 1.081847050930338
$ws.Cells.Item(23, 14).Value = 1.072424412831233
$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.066942010805954
$ws.Cells.Item(24, 4).Value = 1.067873805697268
$ws.Cells.Item(24, 5).Value = 1.071300213131534
$ws.Cells.Item(24, 6).Value = 1.080663226541775
$ws.Cells.Item(24, 9).Value = 1.044219979483893
$ws.Cells.Item(24, 10).Value = 1.07293046756885
$ws.Cells.Item(24, 11).Value = 1.071136281314829
$ws.Cells.Item(24, 12).Value = 1.074551541130275
$ws.Cells.Item(24, 13).Value = 1.083884487062452
$ws.Cells.Item(24, 14).Value = 1.074454151899643
$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.069875852672926
$ws.Cells.Item(25, 4).Value = 1.070190552011498
$ws.Cells.Item(25, 5).Value = 1.073886833690149
$ws.Cells.Item(25, 6).Value = 1.08332796317607
$ws.Cells.Item(25, 9).Value = 1.044897991744726
$ws.Cells.Item(25, 10).Value = 1.075270998762431
$ws.Cells.Item(25, 11).Value = 1.073138784144552
$ws.Cells.Item(25, 12).Value = 1.076824288570984
$ws.Cells.Item(25, 13).Value = 1.086238253158412
$ws.Cells.Item(25, 14).Value = 1.076798006915981
